$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Features")

$ws.Range("A17").Value = 1022
$ws.Range("C17").Value = [double]"3.3688134922222591E-3"
$ws.Range("D17").Value = [double]"9.5667164284161214E-3"
$ws.Range("E17").Value = [double]"5.1865121585082403E-3"
$ws.Range("F17").Value = [double]"1.1102341291114293E-5"
$ws.Range("G17").Value = [double]"4.6030793084526882E-5"
$ws.Range("H17").Value = [double]"2.4555606252025172E-5"

$wb.Save()
